$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=3; A=20150926; E=14},
    @{Row=4; A=20150927; E=15},
    @{Row=5; A=20150928; E=9},
    @{Row=6; A=20150929; E=18},
    @{Row=7; A=20150930; E=5},
    @{Row=8; A=20150931; E=7},
    @{Row=9; A=20150932; E=9},
    @{Row=10; A=20150933; E=19},
    @{Row=11; A=20150934; E=7},
    @{Row=12; A=20150935; E=20},
    @{Row=13; A=20150936; E=14},
    @{Row=14; A=20150937; E=16},
    @{Row=15; A=20150938; E=19},
    @{Row=16; A=20150939; E=10},
    @{Row=17; A=20150940; E=5},
    @{Row=18; A=20150941; E=6},
    @{Row=19; A=20150942; E=19},
    @{Row=20; A=20150943; E=13},
    @{Row=21; A=20150944; E=10},
    @{Row=22; A=20150945; E=7},
    @{Row=23; A=20150946; E=20},
    @{Row=24; A=20150947; E=8},
    @{Row=25; A=20150948; E=8},
    @{Row=26; A=20150949; E=19},
    @{Row=27; A=20150950; E=5},
    @{Row=28; A=20150951; E=14},
    @{Row=29; A=20150952; E=17},
    @{Row=30; A=20150953; E=7},
    @{Row=31; A=20150954; E=10},
    @{Row=32; A=20150955; E=5},
    @{Row=33; A=20150956; E=20},
    @{Row=34; A=20150957; E=5},
    @{Row=35; A=20150958; E=16},
    @{Row=36; A=20150959; E=19},
    @{Row=37; A=20150960; E=9},
    @{Row=38; A=20150961; E=7},
    @{Row=39; A=20150962; E=16},
    @{Row=40; A=20150963; E=14},
    @{Row=41; A=20150964; E=14},
    @{Row=42; A=20150965; E=15},
    @{Row=43; A=20150966; E=15},
    @{Row=44; A=20150967; E=7},
    @{Row=45; A=20150968; E=15},
    @{Row=46; A=20150969; E=12},
    @{Row=47; A=20150970; E=8},
    @{Row=48; A=20150971; E=7},
    @{Row=49; A=20150972; E=16},
    @{Row=50; A=20150973; E=6},
    @{Row=51; A=20150974; E=5},
    @{Row=52; A=20150975; E=17},
    @{Row=53; A=20150976; E=14},
    @{Row=54; A=20150977; E=16},
    @{Row=55; A=20150978; E=6},
    @{Row=56; A=20150979; E=10},
    @{Row=57; A=20150980; E=20},
    @{Row=58; A=20150981; E=12},
    @{Row=59; A=20150982; E=17},
    @{Row=60; A=20150983; E=15},
    @{Row=61; A=20150984; E=16},
    @{Row=62; A=20150985; E=9},
    @{Row=63; A=20150986; E=15}

)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
